$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.81"
$ws.Range("E2").Value = "'2.48%"
$ws.Range("E3").Value = "'4.03%"
$ws.Range("D4").Value = "'5.068"
$ws.Range("E4").Value = "'4.90%"
$ws.Range("D5").Value = "'0.06650"
$ws.Range("E5").Value = "'3.69%"
$ws.Range("D6").Value = "'7.395"
$ws.Range("E6").Value = "'4.71%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.374"
$ws.Range("E7").Value = "'5.80%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9374"
$ws.Range("E8").Value = "'3.62%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1574"
$ws.Range("E9").Value = "'2.25%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.06650"
$ws.Range("E10").Value = "'6.51%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07595"
$ws.Range("E11").Value = "'1.76%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02941"
$ws.Range("E12").Value = "'0.46%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.08987"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001597"
$ws.Range("E14").Value = "'1.10%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04498"
$ws.Range("E15").Value = "'2.22%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006479"
$ws.Range("E16").Value = "'0.66%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006279"
$ws.Range("E17").Value = "'3.13%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.443"
$ws.Range("E18").Value = "'-1.18%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.408"
$ws.Range("E19").Value = "'3.06%"
$ws.Range("D20").Value = "'2.253"
$ws.Range("E20").Value = "'0.85%"
$ws.Range("D21").Value = "'0.3218"
$ws.Range("E21").Value = "'2.33%"
$ws.Range("D22").Value = "'0.1298"
$ws.Range("E22").Value = "'-4.02%"
$ws.Range("D23").Value = "'4.056"
$ws.Range("E23").Value = "'3.51%"
$ws.Range("E24").Value = "'3.21%"
$ws.Range("D25").Value = "'0.001181"
$ws.Range("E25").Value = "'0.43%"
$ws.Range("D26").Value = "'0.004142"
$ws.Range("E26").Value = "'-3.61%"
$ws.Range("D27").Value = "'0.0001247"
$ws.Range("E27").Value = "'5.71%"
$ws.Range("D28").Value = "'0.0001618"
$ws.Range("E28").Value = "'-2.42%"
$ws.Range("D40").Value = "'0.04192"
$ws.Range("E40").Value = "'2.68%"
$ws.Range("D41").Value = "'0.006728"
$ws.Range("E41").Value = "'1.18%"
$ws.Range("D42").Value = "'0.1247"
$ws.Range("E42").Value = "'-11.29%"
$ws.Range("D43").Value = "'0.002016"
$ws.Range("E43").Value = "'-3.56%"
$ws.Range("D44").Value = "'0.01230"
$ws.Range("E44").Value = "'11.24%"
$ws.Range("D45").Value = "'0.00005600"
$ws.Range("E45").Value = "'1.00%"
$ws.Range("E46").Value = "'20.74%"
$ws.Range("D47").Value = "'0.01307"
$ws.Range("E47").Value = "'-29.35%"
